# semana 24 de 2025
# Insert a new data row (evento 100 - Accidente ofidico) right after the
# header row, shifting every existing record down by one row, and update
# the Esperado/Observado/valor p figures for the new reporting week.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row 2, pushing rows 2:33 down to 3:34 -------------------
$ws.Rows("2:2").Insert()

# The insert copies the header row's bold/centered style onto the new row;
# put it back to the plain style used by the rest of the data rows.
$ws.Range("A2:E2").Style = "Normal"

# Column A holds "evento" codes that must stay text (e.g. "100", "113"),
# not be auto-coerced to numbers.
$ws.Range("A2").NumberFormat = "@"

$ws.Range("A2").Value = "100"
$ws.Range("B2").Value = "Accidente ofidico"
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 0

# --- Update Esperado (C) / Observado (D) / valor p (E) for the events -----
# that shifted down (rows now 3..34) wherever the week's figures changed.
$ws.Range("C3").Value = 4
$ws.Range("E3").Value = 0.2

$ws.Range("D4").Value = 1

$ws.Range("C5").Value = 5
$ws.Range("D5").Value = 5
$ws.Range("E5").Value = 0.18

$ws.Range("C6").Value = 1
$ws.Range("D6").Value = 9

$ws.Range("D7").Value = 6
$ws.Range("E7").Value = 0.01

$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 1

$ws.Range("C9").Value = 49
$ws.Range("D9").Value = 35
$ws.Range("E9").Value = 0.01

$ws.Range("D11").Value = 2
$ws.Range("E11").Value = 0.18

$ws.Range("D12").Value = 7
$ws.Range("E12").Value = 0.02

$ws.Range("C13").Value = 28
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0

$ws.Range("C14").Value = 2
$ws.Range("E14").Value = 0.14

$ws.Range("D15").Value = 0

$ws.Range("D16").Value = 4

$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 12
$ws.Range("E17").Value = 0.05

$ws.Range("C18").Value = 3
$ws.Range("E18").Value = 0.05

$ws.Range("C19").Value = 7
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = 0.1

$ws.Range("D20").Value = 1

$ws.Range("D21").Value = 3

$ws.Range("D23").Value = 2
$ws.Range("E23").Value = 0.08

$ws.Range("C24").Value = 0
$ws.Range("D24").Value = 0
$ws.Range("E24").Value = 1

$ws.Range("D28").Value = 0

$ws.Range("C30").Value = 0
$ws.Range("E30").Value = 0

$ws.Range("D32").Value = 7
$ws.Range("E32").Value = 0.15

$ws.Range("C33").Value = 6
$ws.Range("D33").Value = 2
$ws.Range("E33").Value = 0.04

$ws.Range("C34").Value = 9
$ws.Range("D34").Value = 8
$ws.Range("E34").Value = 0.13
